# Update Hotel Details sheet: adjust "Price per Night" values slightly
# and replace "Total Price" values with "N/A" for all three hotels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Holiday Inn NAIROBI TWO RIVERS MALL by IHG
$ws.Range("C2").Value = "₹ 19,986"
$ws.Range("D2").Value = "N/A"

# Row 3 - JW Marriott Hotel Nairobi
$ws.Range("C3").Value = "₹ 203,399"
$ws.Range("D3").Value = "N/A"

# Row 4 - Yaya Hotel & Apartments
$ws.Range("C4").Value = "₹ 19,632"
$ws.Range("D4").Value = "N/A"
